# ------------------------------------------------------------------
# Commit: "xoa cac phan tong cua ti le chiet khau"
# ("remove the total sections of the discount rate")
#
# 1. "Đơn sale chính": Tổng row -> Tỉ lệ chiết khấu sale chính (M4): 0.23 -> 0
# 2. Insert a new sheet "Đơn thu nợ" between "Đơn sale chính" and "Lương"
#    with a debt-collection order ledger (header + 1 data row + Tổng row).
# 3. Rewrite "Lương": drop the "... tại HỆ THỐNG" discount-rate-total block
#    (rows for Chiết khấu/Đơn n bác sĩ/Công phụ phẫu/Ứng lương "tại HỆ THỐNG")
#    and the standalone "Tổng lương tại HỆ THỐNG" row; add new
#    "Chiết khấu thu nợ tại <cơ sở>" rows per branch; refresh totals.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. "Đơn sale chính": Tổng row, Tỉ lệ chiết khấu sale chính (M4): 0.23 -> 0 ---
$wsSale = $wb.Worksheets.Item("Đơn sale chính")
$wsSale.Cells.Item(4, 13).Value = 0

# --- 2. Insert new sheet "Đơn thu nợ" right before "Lương" ---
$wsLuongBefore = $wb.Worksheets.Item("Lương")
$wsDebt = $wb.Worksheets.Add($wsLuongBefore)
$wsDebt.Name = "Đơn thu nợ"

# NOTE: worksheet handles in this host resolve by *position*, not stable object
# identity — after Add() shifts "Lương" from index 2 to index 3, the old
# $wsLuongBefore handle now resolves to whatever sheet sits at index 2 (i.e. the
# freshly added one). Re-fetch "Lương" by name so later writes land on the right
# physical sheet.
$wsLuong = $wb.Worksheets.Item("Lương")

# Row 1: headers (A1:Y1)
$debtHeaders = @(
    "Tiền tố",
    "Mã đơn thu nợ",
    "Lượng thu",
    "Ngày thu",
    "Cơ sở",
    "Đơn nợ",
    "Tên dịch vụ",
    "Khách hàng",
    "Nguồn khách",
    "Sale chính",
    "Đơn giá gốc",
    "Sale phụ",
    "Upsale",
    "Đơn giá",
    "Đã thanh toán",
    "Bác sĩ 1",
    "Bác sĩ 2",
    "Tỉ lệ chiết khấu sale chính",
    "Chiết khấu sale chính",
    "Tỉ lệ chiết khấu sale phụ",
    "Chiết khấu sale phụ",
    "Tỉ lệ chiết khấu bác sĩ 1",
    "Chiết khấu bác sĩ 1",
    "Tỉ lệ chiết khấu bác sĩ 2",
    "Chiết khấu bác sĩ 2"
)
for ($i = 0; $i -lt $debtHeaders.Length; $i++) {
    $wsDebt.Cells.Item(1, $i + 1).Value = $debtHeaders[$i]
}

# Row 2: data row (A2:Y2)
$debtRow2 = @(
    "TN",   # A
    164,   # B
    3000000,   # C
    "'07-25-2024",   # D
    "CẦN THƠ",   # E
    "HD-LUXURY-185",   # F
    "Tiêm Filler",   # G
    "Nguyễn Thị Hồng Trang",   # H
    "Cá nhân",   # I
    "Lê Văn Linh",   # J
    29500000,   # K
    $null,   # L
    $null,   # M
    29500000,   # N
    6000000,   # O
    "CTV Ngoài",   # P
    $null,   # Q
    0.17,   # R
    510000.0000000001,   # S
    0,   # T
    0,   # U
    0,   # V
    0,   # W
    0,   # X
    0   # Y
)
for ($i = 0; $i -lt $debtRow2.Length; $i++) {
    if ($null -ne $debtRow2[$i]) {
        $wsDebt.Cells.Item(2, $i + 1).Value = $debtRow2[$i]
    }
}

# Row 3: "Tổng" summary row (A3:Y3)
$debtRow3 = @(
    "Tổng",   # A
    1,   # B
    3000000,   # C
    $null,   # D
    $null,   # E
    $null,   # F
    $null,   # G
    $null,   # H
    $null,   # I
    $null,   # J
    29500000,   # K
    $null,   # L
    0,   # M
    29500000,   # N
    6000000,   # O
    $null,   # P
    $null,   # Q
    0,   # R
    510000.0000000001,   # S
    0,   # T
    0,   # U
    0,   # V
    0,   # W
    0,   # X
    0   # Y
)
for ($i = 0; $i -lt $debtRow3.Length; $i++) {
    if ($null -ne $debtRow3[$i]) {
        $wsDebt.Cells.Item(3, $i + 1).Value = $debtRow3[$i]
    }
}

# --- 3. Rewrite "Lương": drop the "... tại HỆ THỐNG" discount-rate-total block
#        and the "Tổng lương tại HỆ THỐNG" row; add "Chiết khấu thu nợ tại
#        <cơ sở>" rows per branch; refresh the totals. ---
$wsLuong.Cells.Clear()

$luongLabels = @(
    "Danh mục lương",
    "Ngày công",
    "Phụ cấp",
    "Lương cơ bản tại CẦN THƠ",
    "Chiết khấu sale chính tại CẦN THƠ",
    "Chiết khấu sale phụ tại CẦN THƠ",
    "Đơn 1 bác sĩ tại CẦN THƠ",
    "Đơn 2 bác sĩ tại CẦN THƠ",
    "Công phụ phẫu 1 tại CẦN THƠ",
    "Công phụ phẫu 2 tại CẦN THƠ",
    "Chiết khấu thu nợ tại CẦN THƠ",
    "Ứng lương tại CẦN THƠ",
    "Lương cơ bản tại LONG XUYÊN",
    "Chiết khấu sale chính tại LONG XUYÊN",
    "Chiết khấu sale phụ tại LONG XUYÊN",
    "Đơn 1 bác sĩ tại LONG XUYÊN",
    "Đơn 2 bác sĩ tại LONG XUYÊN",
    "Công phụ phẫu 1 tại LONG XUYÊN",
    "Công phụ phẫu 2 tại LONG XUYÊN",
    "Chiết khấu thu nợ tại LONG XUYÊN",
    "Ứng lương tại LONG XUYÊN",
    "Lương cơ bản tại SÓC TRĂNG",
    "Chiết khấu sale chính tại SÓC TRĂNG",
    "Chiết khấu sale phụ tại SÓC TRĂNG",
    "Đơn 1 bác sĩ tại SÓC TRĂNG",
    "Đơn 2 bác sĩ tại SÓC TRĂNG",
    "Công phụ phẫu 1 tại SÓC TRĂNG",
    "Công phụ phẫu 2 tại SÓC TRĂNG",
    "Chiết khấu thu nợ tại SÓC TRĂNG",
    "Ứng lương tại SÓC TRĂNG",
    "Tổng lương tại CẦN THƠ",
    "Tổng lương tại LONG XUYÊN",
    "Tổng lương tại SÓC TRĂNG",
    "Tổng lương"
)
$luongValues = @(
    7,
    25,
    875000,
    $null,
    1000000,
    0,
    0,
    0,
    0,
    0,
    510000.0000000001,
    -4469000,
    10267857.14285714,
    910000,
    0,
    0,
    0,
    0,
    0,
    0,
    -0.0,
    15401785.71428571,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    -0.0,
    -2084000,
    11177857.14285714,
    15401785.71428571,
    24495642.85714286
)
for ($i = 0; $i -lt $luongLabels.Length; $i++) {
    $wsLuong.Cells.Item($i + 1, 1).Value = $luongLabels[$i]
    if ($null -ne $luongValues[$i]) {
        $wsLuong.Cells.Item($i + 1, 2).Value = $luongValues[$i]
    }
}

